$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.687.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.791.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "

$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.79%  "

$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.047.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.793.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.653.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0531"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.71%  "

$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("E33").Value = "  -2.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.442.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0192"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.640"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.28%  "

$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "85.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.75%  "

$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.96%  "

$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0497"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.945.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.93%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.83%  "
